$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: change the Heading1 text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Superstars Free: Review and Features 2021", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Play Superstars Free | Exciting Features and Unique Design", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2/9. The "Meta description: ..." paragraph (bold "Meta description" run +
#      plain blurb run) right after the title gets removed entirely, and a
#      new paragraph with the *same* bold-run formatting but different text
#      ("Play Superstars Free | Exciting Features and Unique Design") is
#      added further down, right after the last "what we don't like" bullet.
#      We reuse the original paragraph's formatted runs (so the new
#      paragraph keeps the same bold-run / empty-leading-run shape) instead
#      of building it from scratch, then delete the original afterwards.
# ---------------------------------------------------------------------------
$metaParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description:*") {
        $metaParaIndex = $i
        break
    }
}
$metaPara = $d.Paragraphs.Item($metaParaIndex)
$boldHeadingFormat = $metaPara.Range.FormattedText

$lastBulletIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Buy Feature can be pricey*") {
        $lastBulletIndex = $i
        break
    }
}
$lastBullet = $d.Paragraphs.Item($lastBulletIndex)
$lastBullet.Range.InsertParagraphAfter()
$newHeadingPara = $d.Paragraphs.Item($lastBulletIndex + 1)
$newHeadingPara.Style = "Normal"
$newHeadingPara.Range.FormattedText = $boldHeadingFormat

# Retarget the pasted text (still "Meta description: Discover ...") to the
# new heading text, within just that paragraph's range.
$newHeadingPara2 = $d.Paragraphs.Item($lastBulletIndex + 1)
$newHeadingPara2.Range.Find.Execute(
    "Meta description: Discover the game Superstars and its exciting features, including a board game, wheel game, and buy feature. Play for free and win big in 2021.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Superstars Free | Exciting Features and Unique Design", 2) | Out-Null

# Now remove the original "Meta description" paragraph (re-locate it fresh,
# it is still the same paragraph, before the new heading paragraph).
$metaParaIndex2 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description:*") {
        $metaParaIndex2 = $i
        break
    }
}
$d.Paragraphs.Item($metaParaIndex2).Range.Delete()

# ---------------------------------------------------------------------------
# 3-8. "What we like" / "What we don't like" bullet text swaps.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Exciting bonus game and special symbols", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Variety of exciting features", 2) | Out-Null

$d.Content.Find.Execute(
    "Fun and creative features", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Chance to win up to 4,586 times bet", 2) | Out-Null

$d.Content.Find.Execute(
    "High maximum win potential", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Unique combination of characters and icons", 2) | Out-Null

$d.Content.Find.Execute(
    "Diverse and original game design", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Diverse and original design", 2) | Out-Null

$d.Content.Find.Execute(
    "Medium-high volatility may not appeal to all players", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Less frequent prizes", 2) | Out-Null

$d.Content.Find.Execute(
    "Buy Feature can be pricey for some players", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "More risky gameplay", 2) | Out-Null

# ---------------------------------------------------------------------------
# 10. Replace the closing italic "Prompt: ..." paragraph text.
# ---------------------------------------------------------------------------
$oldPrompt = 'Prompt: Create a cartoon-style feature image for the game Superstars that features a happy Maya warrior with glasses. The image should be eye-catching and colorful, highlighting the fun theme of the game. The warrior should be surrounded by the game''s symbols, such as Lady Pig, Finn, Brute, Gonzo, and the eight-sided star. The image should also include the slot''s title, "Superstars," in bold and colorful letters. Make sure the image captures the exciting and creative elements of the game, such as the Bonus Game, Wheel Game, and Replay Feature.'
$newPrompt = "Read our review of Superstars, a slot game with exciting features and a unique design. Play for free now!"
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2) | Out-Null
